$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D (Tipo) to make room for MAE
$ws.Columns.Item(4).Insert()

# Header for new column
$ws.Range("D1").Value = "MAE"

# New MAE values for rows 2-4
$ws.Range("D2").Value = 0.1292627146720763
$ws.Range("D3").Value = 0.1526989685211047
$ws.Range("D4").Value = 0.1401297014068613
